$wb = $excel.ActiveWorkbook

# zh-cn sheet: row 2 (584b702f... zh-cn.xlf) handoff/handback timestamps
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-18 12:48:05"
$wsZhCn.Range("H2").Value = "2016-03-18 12:48:27"

# de-de sheet: row 2 (584b702f... de-de.xlf) handoff/handback timestamps
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-18 12:48:08"
$wsDeDe.Range("H2").Value = "2016-03-18 12:48:32"
